$d = $word.ActiveDocument
$bullet = [char]0x2022

# ---------------------------------------------------------------------
# 1. Collapse the three long "CORE COMPETENCIES" paragraphs into a single
#    short summary paragraph.
# ---------------------------------------------------------------------
$oldCore1 = "Product Management & Strategy: Product Conception & Ideation: Market Research, User Needs Analysis, Opportunity Assessment $bullet Product Architecture & Design: Technical Architecture, System Design, Scalability Planning $bullet Product Lifecycle Management: Roadmap Development, Feature Prioritization, Release Planning $bullet B2B SaaS Development: Multi-tenant Architecture, Platform Strategy, API Design"
$newCore1 = "Product Management & Strategy $bullet Technical Product Development $bullet Platform & Infrastructure"

$range = $d.Content
$range.Find.Execute($oldCore1, $false, $false, $false, $false, $false, $true, 1, $false, $newCore1, 2) | Out-Null

# Delete the two paragraphs that followed (now fully redundant).
$oldCore2 = "Technical Product Development: Full-Stack Development: Python (Django/GeoDjango), JavaScript, React, PostgreSQL $bullet Cloud Platforms: AWS (EC2, RDS, S3), Google Cloud Platform, Microsoft Azure $bullet Big Data Technologies: Apache Spark, PySpark, Hadoop, Snowflake, dbt $bullet API Development: RESTful APIs, Microservices, Third-party Integrations"
$oldCore3 = "Platform & Infrastructure: Multi-tenant Architecture: Scalable Platform Design, Data Isolation, Performance Optimization $bullet Data Warehousing: ETL/ELT Pipelines, Data Governance, Quality Control $bullet Geospatial Platforms: PostGIS, ArcGIS, Quantum GIS, OSGeo Integration $bullet Security & Compliance: Data Privacy, Access Control, Regulatory Compliance"

$findRange = $d.Content
$findRange.Find.Execute($oldCore2) | Out-Null
if ($findRange.Find.Found) {
    $para = $findRange.Paragraphs(1)
    $para.Range.Delete()
}

$findRange2 = $d.Content
$findRange2.Find.Execute($oldCore3) | Out-Null
if ($findRange2.Find.Found) {
    $para2 = $findRange2.Paragraphs(1)
    $para2.Range.Delete()
}

# ---------------------------------------------------------------------
# 2. Insert a new "TECHNICAL SKILLS" section before the closing
#    "For a more detailed..." paragraph.
#    Insert all four paragraphs first (while they still inherit the
#    "Normal" style of the bullet paragraph above them), and only switch
#    the heading paragraph's style afterwards - this keeps the three
#    body paragraphs free of any stray style assignment.
# ---------------------------------------------------------------------
$anchorText = "Trained analytical and engineering staff on open source geospatial technology for analysis, segmentation, and visualization"
$anchorRange = $d.Content
$anchorRange.Find.Execute($anchorText) | Out-Null
$anchorPara = $anchorRange.Paragraphs(1)

$anchorPara.Range.InsertParagraphAfter()
$headingPara = $anchorPara.Next()
$headingPara.Range.Text = "TECHNICAL SKILLS"

$headingPara.Range.InsertParagraphAfter()
$line1 = $headingPara.Next()
$line1.Range.Text = "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development"

$line1.Range.InsertParagraphAfter()
$line2 = $line1.Next()
$line2.Range.Text = "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; API Development"

$line2.Range.InsertParagraphAfter()
$line3 = $line2.Next()
$line3.Range.Text = "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Security & Compliance"

# Now promote the heading paragraph to Heading2 (after all siblings exist).
$headingPara.Style = "Heading2"

Write-Output "Done"
